# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型"
# sheets to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 170
$ws1.Range("F3").Value = 653
$ws1.Range("F5").Value = 220
$ws1.Range("F6").Value = 1546
$ws1.Range("F7").Value = 36
$ws1.Range("F8").Value = 3128
$ws1.Range("F9").Value = 452
$ws1.Range("F10").Value = 724

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 170
$ws4.Range("F3").Value = 653
$ws4.Range("F6").Value = 220
$ws4.Range("F7").Value = 1546
$ws4.Range("F8").Value = 36
$ws4.Range("F9").Value = 3128
$ws4.Range("F10").Value = 452
$ws4.Range("F11").Value = 724

$wb.Save()
